$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(464, 1).Value = "2023-11-25 00:55:53"
$ws.Cells.Item(464, 2).Value = 21
$ws.Cells.Item(464, 3).Value = 5
$ws.Cells.Item(464, 4).Value = 12
$ws.Cells.Item(464, 5).Value = 0
$ws.Cells.Item(464, 6).Value = 3
$ws.Cells.Item(464, 7).Value = 2
$ws.Cells.Item(464, 8).Value = 23
$ws.Cells.Item(464, 9).Value = 20
$ws.Cells.Item(464, 10).Value = 15
$ws.Cells.Item(464, 11).Value = 0.001
$ws.Cells.Item(464, 12).Value = 0.01
$ws.Cells.Item(464, 13).Value = 0.003
$ws.Cells.Item(464, 14).Value = 100
$ws.Cells.Item(464, 15).Value = 512
$ws.Cells.Item(464, 16).Value = 10
$ws.Cells.Item(464, 17).Value = 7
$ws.Cells.Item(464, 18).Value = 0.2380952380952381

$ws.Cells.Item(465, 1).Value = "2023-11-25 11:06:27"
$ws.Cells.Item(465, 2).Value = 56
$ws.Cells.Item(465, 3).Value = 7
$ws.Cells.Item(465, 4).Value = 21
$ws.Cells.Item(465, 5).Value = 1
$ws.Cells.Item(465, 6).Value = 5
$ws.Cells.Item(465, 7).Value = 1
$ws.Cells.Item(465, 8).Value = 41
$ws.Cells.Item(465, 9).Value = 21
$ws.Cells.Item(465, 10).Value = 36
$ws.Cells.Item(465, 11).Value = 0.001
$ws.Cells.Item(465, 12).Value = 0.01
$ws.Cells.Item(465, 13).Value = 0.003
$ws.Cells.Item(465, 14).Value = 100
$ws.Cells.Item(465, 15).Value = 512
$ws.Cells.Item(465, 16).Value = 10
$ws.Cells.Item(465, 17).Value = 7
$ws.Cells.Item(465, 18).Value = 0.125

$ws.Cells.Item(466, 1).Value = "2023-11-25 11:44:58"
$ws.Cells.Item(466, 2).Value = 98
$ws.Cells.Item(466, 3).Value = 13
$ws.Cells.Item(466, 4).Value = 23
$ws.Cells.Item(466, 5).Value = 3
$ws.Cells.Item(466, 6).Value = 3
$ws.Cells.Item(466, 7).Value = 7
$ws.Cells.Item(466, 8).Value = 44
$ws.Cells.Item(466, 9).Value = 58
$ws.Cells.Item(466, 10).Value = 23
$ws.Cells.Item(466, 11).Value = 0.001
$ws.Cells.Item(466, 12).Value = 0.01
$ws.Cells.Item(466, 13).Value = 0.003
$ws.Cells.Item(466, 14).Value = 100
$ws.Cells.Item(466, 15).Value = 512
$ws.Cells.Item(466, 16).Value = 10
$ws.Cells.Item(466, 17).Value = 7
$ws.Cells.Item(466, 18).Value = 0.1326530612244898

$ws.Cells.Item(467, 1).Value = "2023-11-25 19:07:44"
$ws.Cells.Item(467, 2).Value = 14
$ws.Cells.Item(467, 3).Value = 0
$ws.Cells.Item(467, 4).Value = 14
$ws.Cells.Item(467, 5).Value = 0
$ws.Cells.Item(467, 6).Value = 0
$ws.Cells.Item(467, 7).Value = 0
$ws.Cells.Item(467, 8).Value = 14
$ws.Cells.Item(467, 9).Value = 14
$ws.Cells.Item(467, 10).Value = 14
$ws.Cells.Item(467, 11).Value = 0.001
$ws.Cells.Item(467, 12).Value = 0.01
$ws.Cells.Item(467, 13).Value = 0.003
$ws.Cells.Item(467, 14).Value = 100
$ws.Cells.Item(467, 15).Value = 512
$ws.Cells.Item(467, 16).Value = 10
$ws.Cells.Item(467, 17).Value = 7
$ws.Cells.Item(467, 18).Value = 0

$ws.Cells.Item(468, 1).Value = "2023-11-25 19:19:02"
$ws.Cells.Item(468, 2).Value = 63
$ws.Cells.Item(468, 3).Value = 7
$ws.Cells.Item(468, 4).Value = 24
$ws.Cells.Item(468, 5).Value = 1
$ws.Cells.Item(468, 6).Value = 3
$ws.Cells.Item(468, 7).Value = 3
$ws.Cells.Item(468, 8).Value = 52
$ws.Cells.Item(468, 9).Value = 46
$ws.Cells.Item(468, 10).Value = 24
$ws.Cells.Item(468, 11).Value = 0.001
$ws.Cells.Item(468, 12).Value = 0.01
$ws.Cells.Item(468, 13).Value = 0.003
$ws.Cells.Item(468, 14).Value = 100
$ws.Cells.Item(468, 15).Value = 512
$ws.Cells.Item(468, 16).Value = 10
$ws.Cells.Item(468, 17).Value = 7
$ws.Cells.Item(468, 18).Value = 0.1111111111111111

$ws.Cells.Item(469, 1).Value = "2023-11-25 20:35:36"
$ws.Cells.Item(469, 2).Value = 31
$ws.Cells.Item(469, 3).Value = 5
$ws.Cells.Item(469, 4).Value = 8
$ws.Cells.Item(469, 5).Value = 1
$ws.Cells.Item(469, 6).Value = 4
$ws.Cells.Item(469, 7).Value = 0
$ws.Cells.Item(469, 8).Value = 24
$ws.Cells.Item(469, 9).Value = 12
$ws.Cells.Item(469, 10).Value = 31
$ws.Cells.Item(469, 11).Value = 0.001
$ws.Cells.Item(469, 12).Value = 0.01
$ws.Cells.Item(469, 13).Value = 0.003
$ws.Cells.Item(469, 14).Value = 100
$ws.Cells.Item(469, 15).Value = 512
$ws.Cells.Item(469, 16).Value = 10
$ws.Cells.Item(469, 17).Value = 7
$ws.Cells.Item(469, 18).Value = 0.1612903225806452
